$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.852.06"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "1.768.67"
$ws.Range("E3").Value = "  -2.30%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").Value = "'339.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("D7").Value = "'0.3778"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.69%  "
$ws.Range("D8").Value = "'0.3369"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.33%  "
$ws.Range("D9").Value = "'45.59"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.13%  "
$ws.Range("D10").Value = "'1.142"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.75%  "
$ws.Range("D11").Value = "'0.07292"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.66%  "
$ws.Range("E12").Value = "  +4.68%  "
$ws.Range("D13").Value = "'1.002"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").Value = "'6.278"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.65%  "
$ws.Range("D15").Value = "'7.272"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").Value = "1.769.11"
$ws.Range("E16").Value = "  -2.22%  "
$ws.Range("D17").Value = "'0.00001059"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.14%  "
$ws.Range("D18").Value = "'0.06621"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").Value = "'81.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.55%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("D21").Value = "'17.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.45%  "
$ws.Range("D22").Value = "'6.365"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.11%  "
$ws.Range("D23").Value = "27.842.06"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").Value = "'11.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.87%  "
$ws.Range("D25").Value = "'2.390"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("D26").Value = "'1.510"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.42%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'152.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'20.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.80%  "
$ws.Range("D29").Value = "'2.365"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.26%  "
$ws.Range("D30").Value = "1.967.93"
$ws.Range("E30").Value = "  -2.41%  "
$ws.Range("D31").Value = "'133.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.75%  "
$ws.Range("D32").Value = "'4.037"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").Value = "'5.946"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.67%  "
$ws.Range("D34").Value = "'0.08781"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("D35").Value = "'12.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.00%  "
$ws.Range("D36").Value = "'0.02375"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.42%  "
$ws.Range("D37").Value = "'0.6743"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.35%  "
$ws.Range("D38").Value = "'0.06295"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.49%  "
$ws.Range("D39").Value = "'5.220"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.66%  "
$ws.Range("D40").Value = "'0.2130"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.26%  "
$ws.Range("D41").Value = "'1.223"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.75%  "
$ws.Range("D42").Value = "'1.477"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.34%  "
$ws.Range("D43").Value = "'8.137"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.98%  "
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("D45").Value = "'13.89"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.19%  "
$ws.Range("D46").Value = "'0.6143"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.05%  "
$ws.Range("D47").Value = "'3.848"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("D48").Value = "'132.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("D49").Value = "'2.038"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.66%  "
$ws.Range("D50").Value = "'0.07285"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("D51").Value = "'1.191"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.68%  "
